# add-new-employee.xlsx : add new fields (Designation Name, Department Name, Branch)
# - "Designation Code" (col M) -> "Designation Name"
# - "Department Code"  (col N) -> "Department Name"
# - new column "Branch" inserted before the existing "Aadhaar" column (old AC, now AD)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the two existing headers in place (M1/N1).
$ws.Range("M1").Value = "Designation Name"
$ws.Range("N1").Value = "Department Name"

# Insert a brand-new column before AC (old AC = "Aadhaar" shifts to AD).
$ws.Columns("AC").Insert()

# Give it a header and a sensible width (matches the neighbouring "Bank" column).
$ws.Range("AC1").Value = "Branch"
$ws.Columns("AC").ColumnWidth = 20.75

# Reflect the author's final view/selection state (scrolled right, AC1 active).
$ws.Range("R1").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 18
$ws.Range("AC1").Select() | Out-Null
